# Generate Report for Handoff
#
# Refreshes the existing handoff entry (row 2 on every sheet, previously
# describing "1f05dc6c-ec38-449a-89f1-5d2a15e41168.md") so that it now
# describes "69244d70-d9c1-413e-8a07-ce8341d40820.md" with refreshed
# timestamps/xliff hashes, and appends a brand new handoff entry (row 3)
# for "ffff9f3cada0-38e5-41fc-8b96-61c540ca5645.md" on every sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------

# Writes a literal text value into a cell, guaranteeing it stays a text
# (shared-string) cell even for values Excel would normally auto-convert
# (e.g. "True"/"False"/""), and resets the cell style back to the plain
# default style (quote-prefix formatting is not desired here).
function Set-TextValue($range, [string]$val) {
    $range.Value = "'" + $val
    $range.Style = "Normal"
}

# Writes a date/time value as literal text and applies the same
# numeric format used elsewhere in the workbook for date columns.
function Set-DateValue($range, [string]$val) {
    $range.Value = $val
    $range.NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

$oldGuid = "1f05dc6c-ec38-449a-89f1-5d2a15e41168"
$newGuid = "69244d70-d9c1-413e-8a07-ce8341d40820"
$addGuid = "ffff9f3cada0-38e5-41fc-8b96-61c540ca5645"

$commitPrefix = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f28112d136b818c41ce6a0992918519c188c3387/e2e/"

$newRowDate   = "2016-09-06 15:18:32"
$newXliffDate = "2016-09-06 15:18:13"
$zeroDate     = "0001-01-01 00:00:00"

$zhHash = "4716bce0137cd261105be31e9068c9e3ad236ba6"
$deHash = "4716bce0137cd261105be31e9068c9e3ad236ba6"

$zhXlf = "$newGuid.$zhHash.zh-cn.xlf"
$deXlf = "$newGuid.$deHash.de-de.xlf"

# ---------------------------------------------------------------------
# Sheet "Overview" (sheet1) - columns A..G
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# refresh row 2 (rename the file described there)
$wsOverview.Cells.Item(2, 1).Value = "$newGuid.md"
Set-DateValue $wsOverview.Cells.Item(2, 7) $newRowDate

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "$commitPrefix$newGuid.md", "", "", "e2e\$newGuid.md")

# new row 3 for the newly handed-off file
$wsOverview.Cells.Item(3, 1).Value = ".md".Insert(0, "") # placeholder, overwritten below
$wsOverview.Cells.Item(3, 1).Value = "$addGuid.md"
$wsOverview.Cells.Item(3, 3).Value = ".md"
Set-TextValue $wsOverview.Cells.Item(3, 4) ""
$wsOverview.Cells.Item(3, 5).Value = "Ready for handoff"
$wsOverview.Cells.Item(3, 6).Value = "Ready for handoff"
Set-DateValue $wsOverview.Cells.Item(3, 7) $newRowDate

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "$commitPrefix$addGuid.md", "", "", "e2e\$addGuid.md")

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------------
# Sheet "zh-cn" (sheet2) - columns A..P
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# refresh row 2
$wsZh.Cells.Item(2, 1).Value = "$newGuid.md"
$wsZh.Cells.Item(2, 7).Value = $zhXlf
Set-DateValue $wsZh.Cells.Item(2, 8) $newXliffDate
Set-DateValue $wsZh.Cells.Item(2, 11) $zeroDate
Set-TextValue $wsZh.Cells.Item(2, 13) "True"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "$commitPrefix$newGuid.md", "", "", "$newGuid.md")

# new row 3
$wsZh.Cells.Item(3, 2).Value = ".md"
$wsZh.Cells.Item(3, 3).Value = "Ready for handoff"
$wsZh.Cells.Item(3, 4).Value = "e2e"
$wsZh.Cells.Item(3, 5).Value = "ht"
Set-TextValue $wsZh.Cells.Item(3, 6) "True"
$wsZh.Cells.Item(3, 7).Value = $zhXlf
Set-DateValue $wsZh.Cells.Item(3, 8) $newXliffDate
Set-TextValue $wsZh.Cells.Item(3, 9) ""
Set-TextValue $wsZh.Cells.Item(3, 10) ""
Set-DateValue $wsZh.Cells.Item(3, 11) $zeroDate
Set-TextValue $wsZh.Cells.Item(3, 12) ""
Set-TextValue $wsZh.Cells.Item(3, 13) "True"
Set-TextValue $wsZh.Cells.Item(3, 14) ""
Set-TextValue $wsZh.Cells.Item(3, 15) "False"
Set-TextValue $wsZh.Cells.Item(3, 16) ""

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "$commitPrefix$addGuid.md", "", "", "$addGuid.md")

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P3"))

# ---------------------------------------------------------------------
# Sheet "de-de" (sheet3) - columns A..P
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# refresh row 2
$wsDe.Cells.Item(2, 1).Value = "$newGuid.md"
$wsDe.Cells.Item(2, 7).Value = $deXlf
Set-DateValue $wsDe.Cells.Item(2, 11) $zeroDate
Set-TextValue $wsDe.Cells.Item(2, 13) "True"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "$commitPrefix$newGuid.md", "", "", "$newGuid.md")

# new row 3
$wsDe.Cells.Item(3, 2).Value = ".md"
$wsDe.Cells.Item(3, 3).Value = "Ready for handoff"
$wsDe.Cells.Item(3, 4).Value = "e2e"
$wsDe.Cells.Item(3, 5).Value = "ht"
Set-TextValue $wsDe.Cells.Item(3, 6) "True"
$wsDe.Cells.Item(3, 7).Value = $deXlf
Set-DateValue $wsDe.Cells.Item(3, 8) $newRowDate
Set-TextValue $wsDe.Cells.Item(3, 9) ""
Set-TextValue $wsDe.Cells.Item(3, 10) ""
Set-DateValue $wsDe.Cells.Item(3, 11) $zeroDate
Set-TextValue $wsDe.Cells.Item(3, 12) ""
Set-TextValue $wsDe.Cells.Item(3, 13) "True"
Set-TextValue $wsDe.Cells.Item(3, 14) ""
Set-TextValue $wsDe.Cells.Item(3, 15) "False"
Set-TextValue $wsDe.Cells.Item(3, 16) ""

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "$commitPrefix$addGuid.md", "", "", "$addGuid.md")

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P3"))
